$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Töttöröö 31.3.2017 @ 14:50"
$ws.Range("A3").Value = "En ole tyytyväinen. 31.3.2017 @ 14:53"
